$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C5").Value = "Deployment, hosting, webes fizetés"
$ws.Range("C6").Value = "Webes fizetés"
$ws.Range("C7").Value = "Webes fizetés production-ben, Tesztek készítése, Android kliens - architektura megtervezés, app skeleton létrehozása"

$ws.Range("C9").Select()
